$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-21: previously generic "SLR -- " placeholder rows.
# Replace with specific SLR topics (column D first, then column E identifiers)
# to mirror the author's original edit / shared-string ordering.
$ws.Range("D16").Value = "SLR -- Foundations"
$ws.Range("D17").Value = "SLR -- Inference"
$ws.Range("D18").Value = "SLR -- Models"
$ws.Range("D19").Value = "SLR -- Assumptions"
$ws.Range("D21").Value = "SLR -- Summary"
$ws.Range("D20").Value = "SLR -- Transformations"

$ws.Range("E16").Value = "SLRFoundations"
$ws.Range("E17").Value = "SLRInference"
$ws.Range("E18").Value = "SLRModels"
$ws.Range("E19").Value = "SLRAssumptions"
$ws.Range("E20").Value = "SLRTransformations"
$ws.Range("E21").Value = "SLRSummary"

# Move the selection to D22 (matches author's final cursor position).
$ws.Range("D22").Select()
